$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Maria Pires -> João Miguel Lima
$ws.Range("A2").Value = 50017
$ws.Range("B2").Value = "João Miguel Lima"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45085
$ws.Range("G2").Value = 6586.44

# Row 3: Antônio Pimenta -> Carlos Eduardo da Paz
$ws.Range("A3").Value = 65402
$ws.Range("B3").Value = "Carlos Eduardo da Paz"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45093
$ws.Range("G3").Value = 6641.84

# Row 4: Srta. Eduarda Azevedo -> Raul Rodrigues
$ws.Range("A4").Value = 85565
$ws.Range("B4").Value = "Raul Rodrigues"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 5729.25

# Row 5: Sophie Campos -> Bruna Freitas
$ws.Range("A5").Value = 49024
$ws.Range("B5").Value = "Bruna Freitas"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45099
$ws.Range("G5").Value = 5536.97

# Row 6: Ana da Mata -> Lívia da Mota
$ws.Range("A6").Value = 65966
$ws.Range("B6").Value = "Lívia da Mota"
$ws.Range("C6").Value = "P&D"
$ws.Range("F6").Value = 45085
$ws.Range("G6").Value = 6851.94

# Row 7: Kaique Costa -> Arthur Pereira
$ws.Range("A7").Value = 42662
$ws.Range("B7").Value = "Arthur Pereira"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45091
$ws.Range("G7").Value = 7893.09

# Row 8: João Miguel Guerra -> Maria Helena Pacheco
$ws.Range("A8").Value = 83087
$ws.Range("B8").Value = "Maria Helena Pacheco"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Doenca"
$ws.Range("F8").Value = 45079
$ws.Range("G8").Value = 8496.99

# Row 9: Eduarda Borges -> Dom Andrade
$ws.Range("A9").Value = 63181
$ws.Range("B9").Value = "Dom Andrade"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45084
$ws.Range("G9").Value = 6036.26

# Row 10: Ana Cecília Pastor -> Dra. Mariah Alves
$ws.Range("A10").Value = 20255
$ws.Range("B10").Value = "Dra. Mariah Alves"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45093
$ws.Range("G10").Value = 9754.06

# Row 11: Lara Albuquerque -> Beatriz Farias
$ws.Range("A11").Value = 4659
$ws.Range("B11").Value = "Beatriz Farias"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 7407.27
